$p = $ppt.ActivePresentation

# --- Slide 27: "Y-PHYS" TextBox -> "X-PHYS" ---
$s27 = $p.Slides.Item(27)
$sh27 = $s27.Shapes.Item("TextBox 22")
$sh27.TextFrame.TextRange.Text = "X-PHYS"
$sh27.Left = 269.2801207401575
$sh27.Width = 420.4664916929134

# --- Slide 28: "Y-DHAL" TextBox -> "X-DHAL" ---
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item("TextBox 25")
$sh28.TextFrame.TextRange.Text = "X-DHAL"
$sh28.Left = 272.24633771259846
$sh28.Width = 414.534133488189

# --- Slide 30: "X-PHYS-2" TextBox -> "X-DHAL-2" ---
$s30 = $p.Slides.Item(30)
$sh30 = $s30.Shapes.Item("TextBox 30")
$sh30.TextFrame.TextRange.Text = "X-DHAL-2"
$sh30.Left = 223.7145231889764
$sh30.Width = 511.59783927559056
